# SwaadSutra Daily report — fill in today's order data.
#
# Sheet "Daily Orders": header row + the single new order row.
# Sheet "Summary": roll up totals for the new order (1 new order, 105 revenue).
# Sheet "Items Breakdown": per-item quantity/revenue breakdown for the order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Daily Orders
# ---------------------------------------------------------------------------
$orders = $wb.Worksheets.Item("Daily Orders")

$headers1 = @("Order ID", "Date", "Customer", "Flat No", "Phone", "Items", "Total", "Status", "Payment", "Collection Date", "Collection Time", "Notes", "Cancel Reason", "Feedback")
for ($i = 0; $i -lt $headers1.Length; $i++) {
    $orders.Cells.Item(1, $i + 1).Value = $headers1[$i]
}

# Numeric cells
$orders.Range("A2").Value = 15
$orders.Range("G2").Value = 105

# Plain text cells (safe as-is — Excel won't reinterpret these as numbers/dates)
$orders.Range("B2").Value = "2026-01-19 05:39"
$orders.Range("C2").Value = "Prajakta Patil"
$orders.Range("D2").Value = "A 804"
$orders.Range("F2").Value = "Wheat Chapati x5, 1 Plate Bhaji x1"
$orders.Range("H2").Value = "NEW"
$orders.Range("I2").Value = "PENDING"
$orders.Range("K2").Value = "08:00"
$orders.Range("L2").Value = "Less Spicy"

# Text cells that LOOK like numbers/dates — force text format first so Excel
# keeps them as strings instead of silently coercing to a number/date serial.
$orders.Range("E2").NumberFormat = "@"
$orders.Range("E2").Value = "779868817"

$orders.Range("J2").NumberFormat = "@"
$orders.Range("J2").Value = "2026-01-20"

# Empty string cells (Cancel Reason / Feedback — not applicable for a new order)
$emptyStringFormula = '=""'
$orders.Range("M2").Formula = $emptyStringFormula
$orders.Range("N2").Formula = $emptyStringFormula

# ---------------------------------------------------------------------------
# Sheet 2: Summary — roll the new order into the day's totals
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A2").Value = 1    # Total Orders
$summary.Range("B2").Value = 1    # New
$summary.Range("G2").Value = 105  # Total Revenue

# ---------------------------------------------------------------------------
# Sheet 3: Items Breakdown
# ---------------------------------------------------------------------------
$items = $wb.Worksheets.Item("Items Breakdown")

$headers3 = @("Item", "Quantity Ordered", "Revenue")
for ($i = 0; $i -lt $headers3.Length; $i++) {
    $items.Cells.Item(1, $i + 1).Value = $headers3[$i]
}

$items.Range("A2").Value = "Wheat Chapati"
$items.Range("B2").Value = 5
$items.Range("C2").Value = 75

$items.Range("A3").Value = "1 Plate Bhaji"
$items.Range("B3").Value = 1
$items.Range("C3").Value = 30
